# Apply the "MainMenuPageNavigationTest" fixture-data edit to testdata.xlsx.
#
# Summary of the change (from the commit / diff):
#   - A "loadPageVerifyURLAndLogo" sheet already exists.
#   - Two new sheets are added at the end: "pageNavigation" and
#     "pageNavigation (2)", holding page-navigation test data
#     (browser / item name / page title / item title) for the Currys
#     site menu items.
#   - The final active/selected sheet is "pageNavigation".
#   - Each sheet's selection cursor is left on a specific cell, matching
#     where the author happened to click before saving.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) First sheet: just move the selection cursor (it is no longer the
#    active/visible tab once the new sheets are added).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "pageNavigation" sheet - subset with Smart Tech / Home & Outdoor.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "pageNavigation"

$ws2.Range("A1").Value = "browser"
$ws2.Range("B1").Value = "item name"
$ws2.Range("C1").Value = "page title"
$ws2.Range("D1").Value = "item title"

$ws2.Range("A2").Value = "chrome"
$ws2.Range("B2").Value = "Smart Tech"
$ws2.Range("C2").Value = "Smart Tech - Get the latest Smart Tech online here | Currys"
$ws2.Range("D2").Value = "Smart Tech"

$ws2.Range("A3").Value = "chrome"
$ws2.Range("B3").Value = "Home & Outdoor"
$ws2.Range("C3").Value = "Home & Outdoor Accessories | Currys"
$ws2.Range("D3").Value = "home and outdoor"

$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws2.Columns.Item(4).AutoFit() | Out-Null
$ws2.Columns.Item(2).ColumnWidth = 15.26953125
$ws2.Columns.Item(3).ColumnWidth = 53.453125
$ws2.Columns.Item(4).ColumnWidth = 24.7265625
$ws2.Range("C4").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) "pageNavigation (2)" sheet - full list of menu categories.
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet2)
$ws3.Name = "pageNavigation (2)"

$ws3.Range("A1").Value = "browser"
$ws3.Range("B1").Value = "item name"
$ws3.Range("C1").Value = "page title"
$ws3.Range("D1").Value = "item title"

$ws3.Range("A2").Value = "chrome"
$ws3.Range("B2").Value = "Appliances"
$ws3.Range("C2").Value = "Kitchen Appliances, Cookers, Washing Machines | Currys"
$ws3.Range("D2").Value = "Household Appliances"

$ws3.Range("A3").Value = "chrome"
$ws3.Range("B3").Value = "TV & Audio"
$ws3.Range("C3").Value = "TVs, DVD & Blu Ray, Home Cinema & Gaming | Currys"
$ws3.Range("D3").Value = "TV and Home Entertainment"

$ws3.Range("A4").Value = "chrome"
$ws3.Range("B4").Value = "Computing"
$ws3.Range("C4").Value = "Laptops, Tablets, Desktop PCs, Computing Accessories | Currys"
$ws3.Range("D4").Value = "Computing"

$ws3.Range("A5").Value = "chrome"
$ws3.Range("B5").Value = "Gaming"
$ws3.Range("C5").Value = "Gaming | Console and PC Gaming | Currys"
$ws3.Range("D5").Value = "Gaming"

$ws3.Range("A6").Value = "chrome"
$ws3.Range("B6").Value = "Cameras"
$ws3.Range("C6").Value = "Digital Cameras, DSLR, Camcorders, Accessories | Currys"
$ws3.Range("D6").Value = "Cameras and camcorders"

$ws3.Range("A7").Value = "chrome"
$ws3.Range("B7").Value = "Phones"
$ws3.Range("C7").Value = "Mobile Phones & Phones Accessories | Currys Mobile"
$ws3.Range("D7").Value = "Mobile Phones"

$ws3.Range("A8").Value = "chrome"
$ws3.Range("B8").Value = "Smart Tech"
$ws3.Range("C8").Value = "Smart Tech - Get the latest Smart Tech online here | Currys"
$ws3.Range("D8").Value = "Smart Tech"

$ws3.Range("A9").Value = "chrome"
$ws3.Range("B9").Value = "Home & Outdoor"
$ws3.Range("C9").Value = "Home & Outdoor Accessories | Currys"
$ws3.Range("D9").Value = "home and outdoor"

$ws3.Columns.Item(2).AutoFit() | Out-Null
$ws3.Columns.Item(3).AutoFit() | Out-Null
$ws3.Columns.Item(4).AutoFit() | Out-Null
$ws3.Columns.Item(2).ColumnWidth = 15.26953125
$ws3.Columns.Item(3).ColumnWidth = 53.453125
$ws3.Columns.Item(4).ColumnWidth = 24.7265625
$ws3.Range("F17").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) Leave "pageNavigation" as the active tab/sheet, matching the
#    workbook's saved bookViews/activeTab + sheetView/tabSelected state.
# ---------------------------------------------------------------------
$ws2.Activate() | Out-Null
Write-Host "done"
